# Apply the data updates described in the commit:
#   - "tram y te My Hoa" row becomes a duplicate entry for "tram y te My Long"
#     (D5: tramytemyhoa.jpg -> tramytemylong.jpg, F5: tytmyhoa -> tytmylong)
#   - the former "tytmylong" account (F6) is renamed to the new "tytmybinh" station
#
# This mirrors the underlying workbook edits captured in the target XML diff:
#   xl/sharedStrings.xml loses the now-unused "tramytemyhoa.jpg" / "tytmyhoa"
#   strings and gains "tytmybinh"; xl/worksheets/sheet1.xml cells D5, F5, F6
#   are re-pointed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "tramytemylong.jpg"
$ws.Range("F5").Value = "tytmylong"
$ws.Range("F6").Value = "tytmybinh"

# Re-fit column D (Ảnh đại diện) now that it holds a longer/duplicated value,
# matching the author's re-saved column width for that column.
$ws.Columns.Item(4).EntireColumn.AutoFit()

# The author's cursor ended up on F7 when the workbook was re-saved.
$ws.Range("F7").Select() | Out-Null
